$d = $word.ActiveDocument

# 1) SmartCard paragraph - English -> Russian
$old1 = "SmartCard " + [char]0x2014 + " это способ хранения SMART на физической карте. You don" + [char]0x2019 + "t need a phone to hold your SMART. If you have ever had a problem with using a phone inside a building" + [char]0x2026 + "you will get why this is ideal. Only the merchant needs the internet connection. As simple to use as a credit or debit card but not controlled by any 3rd party processor. And not a preloaded card using a visa fiat system" + [char]0x2026 + "this is real crypto to crypto. "
$new1 = "SmartCard " + [char]0x2014 + " это способ хранения SMART на физической карте. Вам не нужен телефон или интернет-соединение. Вероятно, вы сталкивались с проблемой качества мобильного интернета, поэтому вы поймёте наше стремление пойти дальше. Только продавцу нужно иметь интернет-соединение. Это так же просто, как использование дебетовой или кредитной карты, но без третьих сторон или процессоров. Это полностью криптовалютная карта. "

$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) DOWNLOAD SMARTPAY APP -> СКАЧАТЬ SMARTPAY APP
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("DOWNLOAD SMARTPAY APP", $true, $false, $false, $false, $false, $true, 1, $false, "СКАЧАТЬ SMARTPAY APP", 2)

# 3) Accept SmartCash paragraph
$old3 = "Accept SmartCash as a payment option in your business with zero fees using a simple SmartCash app. The SmartPay app is available for use anywhere in the world, all it requires is access to the internet."
$new3 = "Принимайте платежи в SmartCash для вашего бизнеса с практически нулевыми комиссиями, используя простое приложение. Приложение SmartPay доступно для использования во всём мире, вам нужно только интернет-соединение."
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# 4) Available for mobile -> Доступно для мобильного
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute("Available for mobile", $true, $false, $false, $false, $false, $true, 1, $false, "Доступно для мобильного", 2)

# 5) " 3RD PARTY WALLET" (leading non-breaking space) -> "СТОРОННИЕ КОШЕЛЬКИ"
$old5 = [char]0x00A0 + "3RD PARTY WALLET"
$new5 = "СТОРОННИЕ КОШЕЛЬКИ"
$r = $d.Content
$r.Find.ClearFormatting()
$r.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
